$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-10"

# Update the header label in I1 (shared string "2022 (through 09-09)" -> "2022 (through 09-10)")
$ws.Range("I1").Value = "2022 (through 09-10)"

# Update September 2022 count (row 10) and the Total row (row 14) for the 2022 column (I)
$ws.Range("I10").Value = 48
$ws.Range("I14").Value = 1185
